$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 125
$ws1.Range("F6").Value = 750
$ws1.Range("F7").Value = 195
$ws1.Range("F8").Value = 253
$ws1.Range("F9").Value = 1049
$ws1.Range("F11").Value = 346
$ws1.Range("F24").Value = 203
$ws1.Range("F26").Value = 148
$ws1.Range("F31").Value = 1013

# 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 1027
$ws2.Range("F5").Value = 1027
$ws2.Range("F6").Value = 4
$ws2.Range("F17").Value = 965
$ws2.Range("F19").Value = 34
$ws2.Range("F20").Value = 609
$ws2.Range("F26").Value = 2287
$ws2.Range("F29").Value = 17

# 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 2389
$ws3.Range("F6").Value = 979
$ws3.Range("F9").Value = 1228

# 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 2389
$ws4.Range("F8").Value = 979
$ws4.Range("F9").Value = 1228
$ws4.Range("F12").Value = 125
$ws4.Range("F14").Value = 750
$ws4.Range("F15").Value = 195
$ws4.Range("F17").Value = 253
$ws4.Range("F18").Value = 1049
$ws4.Range("F19").Value = 346
$ws4.Range("F21").Value = 1027
$ws4.Range("F31").Value = 203
$ws4.Range("F32").Value = 148
$ws4.Range("F40").Value = 34
$ws4.Range("F46").Value = 1013
